# "Funciones parte 1 corregidas"
# On the "Positions" sheet, insert a new leading index column (A) in the
# 1-6 header/data block: the old A:K columns move right to B:L, the header
# row gets a proper "Type" label (instead of the stray lowercase "type"
# shared-string), and the new column A is filled with a 1..5 row counter
# (blank on the header row). Selection moves from F15 to F11, and column A
# is auto-fitted to its (narrow) content width.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Positions")

# --- 1. snapshot current A:K values for rows 1-6 before moving anything ---
$cols = @("A","B","C","D","E","F","G","H","I","J","K")
$oldVals = @{}
foreach ($r in 1..6) {
    foreach ($col in $cols) {
        $oldVals["$col$r"] = $ws.Range("$col$r").Value2
    }
}

# --- 2. shift columns A..K -> B..L for rows 1-6 (right to left so we never
#        clobber a source cell before it has been read) ---
$mapping = @{
    "A" = "B"; "B" = "C"; "C" = "D"; "D" = "E"; "E" = "F"; "F" = "G";
    "G" = "H"; "H" = "I"; "I" = "J"; "J" = "K"; "K" = "L"
}
for ($i = $cols.Length - 1; $i -ge 0; $i--) {
    $col = $cols[$i]
    $newCol = $mapping[$col]
    for ($r = 1; $r -le 6; $r++) {
        $ws.Range("$newCol$r").Value2 = $oldVals["$col$r"]
    }
}

# --- 3. give the new header cell L1 the same bold header formatting as the
#        rest of row 1 (it had no formatting of its own before the shift) ---
$ws.Range("B1").Copy() | Out-Null
$ws.Range("L1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# --- 4. correct the header text: column E is now "Type" (proper case)
#        instead of the old lowercase "type" label ---
$ws.Range("E1").Value2 = "Type"

# --- 5. clear the now-empty A1 header cell completely (no header above the
#        new index column) ---
$ws.Range("A1").Clear()

# --- 6. fill the new column A with a simple 1..5 row counter ---
for ($r = 2; $r -le 6; $r++) {
    $ws.Range("A$r").Value2 = $r - 1
}

# --- 7. auto-fit the new narrow index column ---
$ws.Columns("A").AutoFit() | Out-Null

# --- 8. move the active selection from F15 to F11 ---
$ws.Range("F11").Select() | Out-Null
